$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.876.02'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.878.63'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.76'
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.04'
$ws.Range("E6").Value = '  +3.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.878.95'
$ws.Range("E7").Value = '  +1.29%  '

$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("E10").Value = '  +3.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.53'
$ws.Range("E11").Value = '  +4.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("E13").Value = '  +16.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.38'
$ws.Range("E14").Value = '  +1.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.527.62'
$ws.Range("E15").Value = '  +1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.902.11'
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.859.63'
$ws.Range("E17").Value = '  +1.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.46'
$ws.Range("E18").Value = '  +2.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.20'
$ws.Range("E21").Value = '  +4.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '474.74'
$ws.Range("E22").Value = '  +1.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.734'
$ws.Range("E23").Value = '  +0.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000164'
$ws.Range("E24").Value = '  +2.75%  '

$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.19'
$ws.Range("E27").Value = '  +0.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.50'
$ws.Range("E28").Value = '  +5.40%  '

$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.025.55'
$ws.Range("E31").Value = '  +1.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.82'
$ws.Range("E32").Value = '  +1.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.61'
$ws.Range("E33").Value = '  +2.45%  '

$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.44'
$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.841.70'
$ws.Range("E36").Value = '  +0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.99'
$ws.Range("E37").Value = '  +23.36%  '

$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("E39").Value = '  +1.67%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.140'
$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.321'
$ws.Range("E43").Value = '  +3.20%  '

$ws.Range("E44").Value = '  +1.18%  '

$ws.Range("E45").Value = '  +13.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '426.26'
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.74'
$ws.Range("E48").Value = '  +2.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.49'
$ws.Range("E49").Value = '  -1.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.61'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("E51").Value = '  +1.60%  '
